# Precinct layer with CT offpcap layer
# Insert a new "major offense per capita" (majoffpc) variable definition row
# into the Variable_Definitions sheet, just above the "borocodenum" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20; this pushes the existing row 20
# ("borocodenum") down to row 21, and inherits formatting from the
# surrounding rows.
$ws.Rows.Item(20).Insert()

# Populate the new row. Setting D20 first, then A20, then B20 reproduces
# the shared-string insertion order seen in the target workbook.
$ws.Range("D20").Value = "major offense per capita"
$ws.Range("A20").Value = "majoffpc"
$ws.Range("B20").Value = "Crime data"
$ws.Range("C20").Value = 0

# Update the active selection to reflect the new layout.
[void]$ws.Range("C15").Select()
